$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the row labels first, in the order that builds up the shared string
# table to match: Current setup, Dac1, I1,uA, DAC2, I2,uA, Divider, a, b, Ur100,mV
$ws.Range("A1").Value = "Current setup"
$ws.Range("A2").Value = "Dac1"
$ws.Range("A4").Value = "I1,uA"
$ws.Range("A5").Value = "DAC2"
$ws.Range("A7").Value = "I2,uA"
$ws.Range("A8").Value = "Divider"
$ws.Range("A9").Value = "a"
$ws.Range("A10").Value = "b"
$ws.Range("A3").Value = "Ur100,mV"
$ws.Range("A6").Value = "Ur100,mV"

# Numeric inputs
$ws.Range("B2").Value = 500
$ws.Range("B3").Value = 53.2
$ws.Range("B5").Value = 3000
$ws.Range("B6").Value = 298.4

# Formulas
$ws.Range("B4").Formula = "=B3*10"
$ws.Range("B7").Formula = "=B6*10"
$ws.Range("B8").Formula = "=B7-B4"
$ws.Range("B9").Formula = "=B5-B2"
$ws.Range("B10").Formula = "=B2*B8-B4*B9"
